$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.799.92"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.308.14"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'323.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").Value = "'105.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "'40.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.37%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "'8.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.31%  "
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "'0.978"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "'15.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "2.658.39"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "2.303.07"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "42.871.29"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "'7.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "'13.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +32.02%  "
$ws.Range("D22").Value = "'73.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "'3.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").Value = "'274.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").Value = "'10.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "'22.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "'38.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.25%  "
$ws.Range("D31").Value = "'166.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").Value = "'6.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.74%  "
$ws.Range("D33").Value = "'0.0890"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "'0.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  -11.85%  "
$ws.Range("D37").Value = "'4.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").Value = "'0.0356"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("D39").Value = "'3.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("D41").Value = "'1.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.48%  "
$ws.Range("D42").Value = "'102.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'70.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  +6.95%  "
$ws.Range("D47").Value = "'83.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.07%  "
$ws.Range("D48").Value = "'114.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "'8.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'5.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").Value = "1.595.12"
$ws.Range("E51").Value = "  +4.93%  "
